$wb = $excel.ActiveWorkbook

# Switch to the "Repayment Schedule" sheet and insert a new blank column
# before column N (shifting old N/O/P -> O/P/Q).
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Select()
$ws.Columns("N:N").Insert()

# Update the selection on the Repayment Schedule sheet.
$ws.Range("R4").Select()
